$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report generation timestamp
$ws.Range("A1").Value = "CreatedAt: 2025-08-23T17:07:24"

# Refreshed intertie pricing figures (columns U-Z) pulled from the latest IESO run
$updates = @{
    "U4" = 33.92
    "V4" = 34.39
    "W4" = 31.26
    "X4" = 34.41
    "Y4" = 32.34
    "Z4" = 34.41
    "U6" = -1.02
    "V6" = -1.17
    "X6" = -1.27
    "Y6" = -1.07
    "Z6" = -1.41
    "U8" = 0
    "V8" = 0
    "U9" = 51
    "V9" = 55.31
    "W9" = 30.79
    "X9" = 34.95
    "Y9" = 33.54
    "Z9" = 35.93
    "U11" = -1.12
    "V11" = -1.27
    "W11" = -1.85
    "Y11" = 0.13
    "Z11" = 0.11
    "U12" = 17.18
    "V12" = 21.02
    "U13" = 0
    "V13" = 0
    "U14" = 51
    "V14" = 55.31
    "W14" = 30.79
    "X14" = 34.95
    "Z14" = 112.69
    "Z15" = 76.75
    "U16" = -1.12
    "V16" = -1.27
    "W16" = -1.85
    "Y16" = 0.13
    "Z16" = 0.11
    "U17" = 17.18
    "V17" = 21.02
    "U18" = 0
    "V18" = 0
    "U19" = 33.89
    "V19" = 34.33
    "W19" = 31.11
    "X19" = 34.28
    "Z19" = 34.35
    "U21" = -1.05
    "V21" = -1.24
    "W21" = -1.52
    "X21" = -1.41
    "Y21" = -1.16
    "Z21" = -1.48
    "U23" = 0
    "V23" = 0
    "U24" = 33.89
    "V24" = 34.33
    "W24" = 31.11
    "X24" = 34.28
    "Z24" = 34.35
    "U26" = -1.05
    "V26" = -1.24
    "W26" = -1.52
    "X26" = -1.41
    "Y26" = -1.16
    "Z26" = -1.48
    "U28" = 0
    "V28" = 0
    "U29" = 33.92
    "V29" = 34.29
    "W29" = 31.02
    "X29" = 34.15
    "Y29" = 32.09
    "Z29" = 34.22
    "U31" = -1.02
    "V31" = -1.27
    "W31" = -1.61
    "X31" = -1.54
    "Y31" = -1.32
    "Z31" = -1.61
    "U33" = 0
    "V33" = 0
    "V34" = 34.33
    "W34" = 30.56
    "X34" = 35.51
    "Y34" = 34.34
    "Z35" = 76.75
    "U36" = -1.15
    "V36" = -1.24
    "W36" = -2.08
    "X36" = -0.18
    "Y36" = 0.93
    "Z36" = 1.15
    "U37" = 26.81
    "U38" = 0
    "V38" = 0
    "U39" = 33.92
    "V39" = 34.39
    "W39" = 31.26
    "X39" = 34.41
    "Y39" = 32.34
    "Z39" = 34.41
    "U41" = -1.02
    "V41" = -1.17
    "X41" = -1.27
    "Y41" = -1.07
    "Z41" = -1.41
    "U43" = 0
    "V43" = 0
    "U44" = 35.08
    "V44" = 35.67
    "W44" = 32.7
    "X44" = 35.76
    "Y44" = 33.44
    "Z44" = 35.65
    "U46" = 0.14
    "V46" = 0.11
    "Y46" = 0.03
    "Z46" = -0.18
    "U48" = 0
    "V48" = 0
    "U49" = 36.06
    "V49" = 36.62
    "W49" = 33.54
    "X49" = 36.49
    "Y49" = 34.37
    "Z49" = 37.08
    "U51" = 1.12
    "V51" = 1.06
    "W51" = 0.91
    "Y51" = 0.96
    "Z51" = 1.26
    "U53" = 0
    "V53" = 0
    "U54" = 33.14
    "V54" = 34.3
    "W54" = 31.69
    "X54" = 35.4
    "Y54" = 33.58
    "Z54" = 36.15
    "U56" = -1.63
    "V56" = -0.93
    "W56" = -0.95
    "Y56" = 0.17
    "Z56" = 0.33
    "U57" = -0.17
    "V57" = -0.32
    "U58" = 0
    "V58" = 0
    "U59" = 36.21
    "V59" = 36.81
    "W59" = 33.86
    "X59" = 36.98
    "Y59" = 34.59
    "Z59" = 36.93
    "U61" = 1.27
    "V61" = 1.25
    "W61" = 1.22
    "Y61" = 1.18
    "Z61" = 1.11
    "U63" = 0
    "V63" = 0
    "U64" = 36.89
    "V64" = 37.47
    "W64" = 34.46
    "X64" = 37.61
    "Y64" = 35.13
    "Z64" = 37.51
    "U66" = 1.96
    "V66" = 1.91
    "W66" = 1.83
    "X66" = 1.92
    "Y66" = 1.72
    "Z66" = 1.69
    "U68" = 0
    "V68" = 0
    "U69" = 37.21
    "V69" = 37.83
    "W69" = 34.79
    "X69" = 38.01
    "Y69" = 35.58
    "Z69" = 37.99
    "U71" = 2.27
    "V71" = 2.27
    "W71" = 2.16
    "X71" = 2.32
    "Y71" = 2.17
    "Z71" = 2.17
    "U73" = 0
    "V73" = 0
    "U74" = 36.09
    "V74" = 36.81
    "W74" = 33.86
    "X74" = 36.94
    "Y74" = 34.59
    "Z74" = 36.93
    "U76" = 1.16
    "V76" = 1.25
    "W76" = 1.22
    "X76" = 1.26
    "Y76" = 1.18
    "Z76" = 1.11
    "U78" = 0
    "V78" = 0
    "U79" = 36.32
    "V79" = 37.02
    "W79" = 34.06
    "X79" = 37.19
    "Y79" = 34.8
    "Z79" = 37.18
    "U81" = 1.38
    "V81" = 1.46
    "W81" = 1.42
    "X81" = 1.51
    "Y81" = 1.39
    "Z81" = 1.36
    "U83" = 0
    "V83" = 0
    "U84" = 32.96
    "V84" = 35.28
    "W84" = 32.35
    "X84" = 36.08
    "Y84" = 34.2
    "Z84" = 36.86
    "U86" = -1.98
    "W86" = -0.29
    "Z86" = 1.03
    "U88" = 0
    "V88" = 0
    "U89" = 33.92
    "V89" = 34.29
    "W89" = 31.02
    "X89" = 34.15
    "Y89" = 32.09
    "Z89" = 34.22
    "U91" = -1.02
    "V91" = -1.27
    "W91" = -1.61
    "X91" = -1.54
    "Y91" = -1.32
    "Z91" = -1.61
    "U93" = 0
    "V93" = 0
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
